$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.619.67'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').Value = '1.958.31'
$ws.Range('E3').Value = '  +0.92%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '244.49'
$ws.Range('E5').Value = '  +0.75%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.616'
$ws.Range('E6').Value = '  +0.99%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '61.67'
$ws.Range('E7').Value = '  +8.05%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  +4.77%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0795'
$ws.Range('E10').Value = '  -6.32%  '
$ws.Range('E11').Value = '  +0.00%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '14.25'
$ws.Range('E12').Value = '  +6.07%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '22.04'
$ws.Range('E13').Value = '  +3.42%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.832'
$ws.Range('E14').Value = '  +2.90%  '
$ws.Range('D15').Value = '2.233.72'
$ws.Range('E15').Value = '  +0.42%  '
$ws.Range('E16').Value = '  +2.99%  '
$ws.Range('D17').Value = '1.959.99'
$ws.Range('E17').Value = '  +1.10%  '
$ws.Range('D18').Value = '36.545.49'
$ws.Range('E18').Value = '  +0.41%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '69.74'
$ws.Range('E19').Value = '  +0.87%  '
$ws.Range('E20').Value = '  -1.11%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '230.22'
$ws.Range('E21').Value = '  +1.23%  '
$ws.Range('E22').Value = '  +1.76%  '
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('E24').Value = '  +4.91%  '
$ws.Range('E25').Value = '  +2.81%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.142'
$ws.Range('E26').Value = '  +7.19%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.18'
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '160.40'
$ws.Range('E28').Value = '  -0.37%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.41'
$ws.Range('E29').Value = '  +0.98%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.30'
$ws.Range('E30').Value = '  +18.30%  '
$ws.Range('E31').Value = '  +1.46%  '
$ws.Range('E32').Value = '  +4.82%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0615'
$ws.Range('E33').Value = '  -0.29%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.45'
$ws.Range('E34').Value = '  +6.98%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.54'
$ws.Range('E35').Value = '  +13.22%  '
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.27'
$ws.Range('E36').Value = '  +4.35%  '
$ws.Range('B37').Value = 'BinanceUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  -0.17%  '
$ws.Range('E38').Value = '  -1.07%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.53'
$ws.Range('E39').Value = '  -9.22%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0980'
$ws.Range('E40').Value = '  -0.82%  '
$ws.Range('E41').Value = '  +0.46%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.17'
$ws.Range('E42').Value = '  +2.10%  '
$ws.Range('E43').Value = '  +0.80%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '16.03'
$ws.Range('E44').Value = '  +2.75%  '
$ws.Range('D45').Value = '1.367.78'
$ws.Range('E45').Value = '  +2.06%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '88.64'
$ws.Range('E46').Value = '  +2.64%  '
$ws.Range('E47').Value = '  +1.06%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '7.13'
$ws.Range('E48').Value = '  +0.16%  '
$ws.Range('E49').Value = '  +0.55%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '45.69'
$ws.Range('E50').Value = '  +5.97%  '
$ws.Range('D51').Value = '2.126.39'
$ws.Range('E51').Value = '  +0.54%  '
